# Insert a new weekly price record for "Femacal de La Calera - Ciboulette"
# above the existing row 108, pushing all subsequent rows (old 108-207) down
# by one (new 109-208). The new row duplicates the data of the row that used
# to be at 108 (now at 109) except for the date (column D) and volume
# (column J), which get fresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 108 (and everything below it) down by one row.
$ws.Rows.Item(108).Insert()

# Seed the new row 108 with the same record that is now sitting at row 109
# (formerly row 108), then tweak the two cells that actually differ.
$ws.Range("A109:R109").Copy()
$ws.Range("A108").PasteSpecial()

$ws.Cells.Item(108, 4).Value = 44512
$ws.Cells.Item(108, 10).Value = 160
